# Regenerate save_data to use K instead of Strike#, update column G (K) values
# per recalculated s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 2
    11 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
